# Update cryptos worksheet with latest prices and volume changes (data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "96.701.36"
$ws.Range("E2").Value = "  +1.59%  "

# Row 3
$ws.Range("D3").Value = "3.588.15"
$ws.Range("E3").Value = "  -0.07%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.21"
$ws.Range("E5").Value = "  +0.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "659.26"
$ws.Range("E6").Value = "  +1.53%  "

# Row 7
$ws.Range("E7").Value = "  +7.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.408"
$ws.Range("E8").Value = "  +0.10%  "

# Row 10
$ws.Range("E10").Value = "  +5.44%  "

# Row 11
$ws.Range("D11").Value = "3.588.30"
$ws.Range("E11").Value = "  +0.06%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.56"
$ws.Range("E12").Value = "  +1.44%  "

# Row 13
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("E14").Value = "  +1.53%  "

# Row 15
$ws.Range("D15").Value = "4.248.94"
$ws.Range("E15").Value = "  -0.45%  "

# Row 16
$ws.Range("D16").Value = "96.524.11"
$ws.Range("E16").Value = "  +1.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000259"
$ws.Range("E17").Value = "  +1.34%  "

# Row 18
$ws.Range("D18").Value = "3.586.98"
$ws.Range("E18").Value = "  -0.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  -2.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  +1.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.90"
$ws.Range("E21").Value = "  -0.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.495"
$ws.Range("E22").Value = "  +2.14%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "513.82"
$ws.Range("E23").Value = "  +0.74%  "

# Row 24
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.45"
$ws.Range("E24").Value = "  -0.48%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000201"
$ws.Range("E25").Value = "  +3.03%  "

# Row 26
$ws.Range("E26").Value = "  +3.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "96.98"
$ws.Range("E27").Value = "  +4.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.81"
$ws.Range("E28").Value = "  +0.12%  "

# Row 29
$ws.Range("D29").Value = "3.778.75"
$ws.Range("E29").Value = "  +0.14%  "

# Row 30
$ws.Range("E30").Value = "  -2.09%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.150"
$ws.Range("E31").Value = "  +8.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.57"

# Row 34
$ws.Range("E34").Value = "  +5.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.16%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "31.78"
$ws.Range("E36").Value = "  +0.38%  "

# Row 37
$ws.Range("E37").Value = "  +1.89%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.54"
$ws.Range("E38").Value = "  +4.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "599.26"
$ws.Range("E39").Value = "  +7.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.61"
$ws.Range("E40").Value = "  +10.30%  "

# Row 41
$ws.Range("E41").Value = "  +0.10%  "

# Row 42
$ws.Range("E42").Value = "  +0.69%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.910"
$ws.Range("E43").Value = "  -1.86%  "

# Row 44
$ws.Range("E44").Value = "  +7.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.78"
$ws.Range("E45").Value = "  +1.79%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "34.61"
$ws.Range("E46").Value = "  +5.61%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.29"
$ws.Range("E47").Value = "  +1.19%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0423"
$ws.Range("E48").Value = "  +1.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.54"
$ws.Range("E49").Value = "  -0.74%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.64"
$ws.Range("E50").Value = "  +5.77%  "

# Row 51
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.28"
$ws.Range("E51").Value = "  +2.70%  "
